$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the History dialogue lines (column B) to the revised wording.
$ws.Range("B4").Value = "It should’ve been during 3-5 PM in the main hall. After serving you tea in the main hall, I didn’t see the Lord again."
$ws.Range("B7").Value = "Near 7 PM, I went to the banquet hall with Butler He."
$ws.Range("B9").Value = "I was with Butler He the entire time."
$ws.Range("B10").Value = "Waiting for the Lord, reheating the dishes in the kitchen——we were always together."
$ws.Range("B11").Value = "It wasn’t until Butler He told us to go to the meeting hall that I found out the Lord had been killed."
$ws.Range("B13").Value = "I’m sorry, I’m just a low-ranking servant, I really don’t know anything."

# Row 7 text got shorter (single line instead of two), so it no longer needs extra height.
$ws.Rows.Item(7).RowHeight = 17

# Restore the active selection to B7, matching the saved view state.
$ws.Range("B7").Select() | Out-Null
